# StockBuyOrderExport.xlsx localisation pass: Russian -> Turkish ("simaris idarecisi")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (товары -> mallar)
$ws.Name = 'mallar'

# Translate the shared header/title strings (Russian -> Turkish)
$ws.Range("A1").Value = 'Tedarikçige Sımarışlar {$v->date}'

$ws.Range("A2").Value = 'İzaat'
$ws.Range("B2").Value = 'Kod'
$ws.Range("C2").Value = 'İsim'
$ws.Range("D2").Value = 'Sımarış'
$ws.Range("E2").Value = 'İhtiyaç'
$ws.Range("F2").Value = 'Fiyat'
$ws.Range("G2").Value = 'Tedarikçi fiyat variantları'

# Row 2 no longer needs the taller custom height
$ws.Rows.Item(2).RowHeight = 13.8

# Move the viewport / selection: topLeftCell A1 -> C1, activeCell F4 -> G11
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G11").Select() | Out-Null
